# coiba_camtrap_ids_gps.xlsx — "redoing sex bias analyses and cleaning script.
# Got until halfway H1" — fill in column H (tool_site) with 0 for the newly
# re-walked rows (94, 96-136, 139) and flip the already-filled H138 back to 0
# (it was prematurely marked 1), then leave the view scrolled/selected where
# the edit left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are missing an H (tool_site) value entirely -> set to 0.
$rowsNeedingH = @(94) + (96..136) + @(139)
foreach ($r in $rowsNeedingH) {
    $ws.Cells.Item($r, 8).Value = 0
}

# Row 138 already had H filled in (as 1) — correct it back to 0.
$ws.Range("H138").Value = 0

# Leave the sheet scrolled/selected on the range the author was working
# through (H95:H142), matching where work left off.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H95:H142").Select()
